$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# Column C (GDP) value updates
$ws.Range("C2").Value = 4043.14643519898
$ws.Range("C3").Value = 2934.187009790061
$ws.Range("C4").Value = 1873.394108966653
$ws.Range("C5").Value = 1460.056109840828
$ws.Range("C6").Value = 1909.084588129339
$ws.Range("C7").Value = 6128.19547247793
$ws.Range("C8").Value = 4547.50930098406
$ws.Range("C9").Value = 4729.735976516416
$ws.Range("C10").Value = 4141.976353364595
$ws.Range("C11").Value = 473.2998774917226
$ws.Range("C12").Value = 1280.225469721551
$ws.Range("C13").Value = 341.5541149051794
$ws.Range("C14").Value = 2983.242707849043
$ws.Range("C15").Value = 665.6274194933962
$ws.Range("C16").Value = 1904.346464968814
$ws.Range("C17").Value = 1503.870423231357
$ws.Range("C18").Value = 1955.461557360978
$ws.Range("C19").Value = 6336.709213679884
$ws.Range("C20").Value = 4633.590358399045
$ws.Range("C21").Value = 5082.354756663512
$ws.Range("C22").Value = 487.7306818514292
$ws.Range("C23").Value = 369.2024078290272
$ws.Range("C24").Value = 0
$ws.Range("C25").Value = 691.8942672110555
$ws.Range("C26").Value = 2024.117324382548
$ws.Range("C27").Value = 11627.81065059172
$ws.Range("C28").Value = 4921.848409120176
$ws.Range("C29").Value = 5360.226632400601
$ws.Range("C30").Value = 2094.024217383061
$ws.Range("C31").Value = 5642.578115155247
$ws.Range("C32").Value = 495.763971160512
$ws.Range("C33").Value = 389.9389667216314
$ws.Range("C34").Value = 5122.180090208862
$ws.Range("C35").Value = 2379.668184479739
$ws.Range("C36").Value = 482.9237812079122
$ws.Range("C37").Value = 2201.396847776877
$ws.Range("C38").Value = 5919.20956823756
$ws.Range("C39").Value = 503.3023574516347
$ws.Range("C40").Value = 419.1838602515346
$ws.Range("C41").Value = 5295.682695961288
$ws.Range("C42").Value = 2497.68592515536
$ws.Range("C43").Value = 493.8183694827482
$ws.Range("C44").Value = 2286.013198234259
$ws.Range("C45").Value = 3252.634165082374
$ws.Range("C46").Value = 449.4203771491282
$ws.Range("C47").Value = 730.3063521039821
$ws.Range("C48").Value = 711.3043470146426
$ws.Range("C49").Value = 515.8271637832048
$ws.Range("C50").Value = 5996.49696468919
$ws.Range("C51").Value = 6103.744960203087
$ws.Range("C52").Value = 2361.056581219794
$ws.Range("C53").Value = 3314.741082534716
$ws.Range("C54").Value = 482.6390663355013
$ws.Range("C55").Value = 729.1196658666737
$ws.Range("C56").Value = 731.9993357350996
$ws.Range("C57").Value = 517.8609592583078
$ws.Range("C58").Value = 6114.227214287786
$ws.Range("C59").Value = 6249.151036691844

# AL column updates (0 -> 1)
$ws.Range("AL15").Value = 1
$ws.Range("AL25").Value = 1
$ws.Range("AL47").Value = 1
$ws.Range("AL55").Value = 1
